# Updated PCM image with new suggestions
# Re-positions/resizes several shapes on slide 2 and tweaks some text
# (capitalisation / wording) on the "Automated Alpha Attribute Selection"
# and "Identifying Comparable Sub-logs" process boxes, plus tightens up
# the "In Depth Comparison" box (dropping its leading blank line and
# hyphenating the title).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1. Move / resize shapes (values chosen in points; PowerPoint stores
#    Left/Top/Width/Height as EMU/12700 internally).
# ---------------------------------------------------------------------

$shInterp = $s.Shapes.Item("Rounded Rectangle 25")
$shInterp.Left   = 189.751047
$shInterp.Top    = 202.889939
$shInterp.Width  = 174.733881
$shInterp.Height = 87.771516

$shAlpha = $s.Shapes.Item("Round Diagonal Corner of Rectangle 3")
$shAlpha.Left   = 461.665137
$shAlpha.Top    = 66.694743
$shAlpha.Width  = 203.104352
$shAlpha.Height = 112.826950

$shDepth = $s.Shapes.Item("Round Diagonal Corner of Rectangle 4")
$shDepth.Left   = 441.714821
$shDepth.Width  = 217.995220
$shDepth.Height = 86.152461

$shSubLogs = $s.Shapes.Item("Round Diagonal Corner of Rectangle 6")
$shSubLogs.Width  = 217.995220
$shSubLogs.Height = 90.685691

$shArrow16 = $s.Shapes.Item("Right Arrow 16")
$shArrow16.Left = 677.389082
$shArrow16.Top  = 241.174588

$shResults = $s.Shapes.Item("Rectangle 1")
$shResults.Left   = 11.087737
$shResults.Top    = 208.430098
$shResults.Width  = 101.433405
$shResults.Height = 82.231361

$shArrow15 = $s.Shapes.Item("Right Arrow 15")
$shArrow15.Left   = 686.968060
$shArrow15.Top    = 113.100258
$shArrow15.Width  = 43.966872
$shArrow15.Height = 20.015927

$shArrow17 = $s.Shapes.Item("Right Arrow 17")
$shArrow17.Left  = 382.163952
$shArrow17.Top   = 241.255767
$shArrow17.Width = 41.871832

$shArrow18 = $s.Shapes.Item("Right Arrow 18")
$shArrow18.Left = 130.200178
$shArrow18.Top  = 241.174588

# ---------------------------------------------------------------------
# 2. Text tweaks.
# ---------------------------------------------------------------------

# "Automated Alpha Attribute Selection" box: re-case the two bullets.
$alphaTextRange = $shAlpha.TextFrame.TextRange
$rankPara = $alphaTextRange.Paragraphs(2)
$alphaTextRange.Characters($rankPara.Start, $rankPara.Length).Text = "Rank trace attributes"
$instPara = $alphaTextRange.Paragraphs(3)
$alphaTextRange.Characters($instPara.Start, $instPara.Length).Text = "Instantiating new alpha attributes"

# "In Depth Comparison" box: drop the leading blank paragraph and hyphenate the title.
$depthTextRange = $shDepth.TextFrame.TextRange
$depthTextRange.Paragraphs(1).Delete()
$titlePara = $depthTextRange.Paragraphs(1)
$depthTextRange.Characters($titlePara.Start, $titlePara.Length).Text = "In-Depth Comparison"

# "Identifying Comparable Sub-logs" box: re-case the title and two bullets.
$subLogsTextRange = $shSubLogs.TextFrame.TextRange
$titlePara2 = $subLogsTextRange.Paragraphs(1)
$subLogsTextRange.Characters($titlePara2.Start, $titlePara2.Length).Text = "Identifying Comparable sub-Logs"
$clusterPara = $subLogsTextRange.Paragraphs(2)
$subLogsTextRange.Characters($clusterPara.Start, $clusterPara.Length).Text = "Stochastic clustering"
$pairPara = $subLogsTextRange.Paragraphs(3)
$subLogsTextRange.Characters($pairPara.Start, $pairPara.Length).Text = "Pair comparison"
